$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the asset code value in B2 (02.01.0001 -> 02.01.0012), preserving the trailing tab
$ws.Range("B2").Value = "02.01.0012`t"

# Move the active selection from K3 to B3
$ws.Range("B3").Select()
